# Update metrics for row 31 (cohort 2022, period_index 3)
$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("C31").Value = 58
$ws.Range("E31").Value = 58 / 2312

# Update metrics for row 37 (cohort 2025, period_index 0)
$ws.Range("C37").Value = 974
$ws.Range("D37").Value = 974
